$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- percentil_dist_25 (column F, index 6): same corrected value for every data row (2-171) ---
$newPercentilDist25 = 0.8762693424044758
for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 6).Value = $newPercentilDist25
}

# --- timestamp (column Z, index 26): one fresh ISO-ish timestamp per row ---
# Row number -> new timestamp string, taken straight from the recorded run.
$rowTimestamps = @{
    2 = "2025-10-19T23:54:27.647356"
    3 = "2025-10-19T23:54:27.647356"
    4 = "2025-10-19T23:54:27.647356"
    5 = "2025-10-19T23:54:27.648355"
    6 = "2025-10-19T23:54:27.648355"
    7 = "2025-10-19T23:54:27.649357"
    8 = "2025-10-19T23:54:27.649357"
    9 = "2025-10-19T23:54:27.649357"
    10 = "2025-10-19T23:54:27.649357"
    11 = "2025-10-19T23:54:27.649357"
    12 = "2025-10-19T23:54:27.650357"
    13 = "2025-10-19T23:54:27.650357"
    14 = "2025-10-19T23:54:27.650357"
    15 = "2025-10-19T23:54:27.650357"
    16 = "2025-10-19T23:54:27.651354"
    17 = "2025-10-19T23:54:27.651354"
    18 = "2025-10-19T23:54:27.651354"
    19 = "2025-10-19T23:54:27.651354"
    20 = "2025-10-19T23:54:27.651354"
    21 = "2025-10-19T23:54:27.652353"
    22 = "2025-10-19T23:54:27.652353"
    23 = "2025-10-19T23:54:27.652353"
    24 = "2025-10-19T23:54:27.652353"
    25 = "2025-10-19T23:54:27.652353"
    26 = "2025-10-19T23:54:27.653352"
    27 = "2025-10-19T23:54:27.653352"
    28 = "2025-10-19T23:54:27.653352"
    29 = "2025-10-19T23:54:27.653352"
    30 = "2025-10-19T23:54:27.654352"
    31 = "2025-10-19T23:54:27.654352"
    32 = "2025-10-19T23:54:27.654352"
    33 = "2025-10-19T23:54:27.654352"
    34 = "2025-10-19T23:54:27.654352"
    35 = "2025-10-19T23:54:27.655350"
    36 = "2025-10-19T23:54:27.656371"
    37 = "2025-10-19T23:54:27.656371"
    38 = "2025-10-19T23:54:27.657359"
    39 = "2025-10-19T23:54:27.657359"
    40 = "2025-10-19T23:54:27.657359"
    41 = "2025-10-19T23:54:27.657359"
    42 = "2025-10-19T23:54:27.657359"
    43 = "2025-10-19T23:54:27.658355"
    44 = "2025-10-19T23:54:27.658355"
    45 = "2025-10-19T23:54:27.658355"
    46 = "2025-10-19T23:54:27.658355"
    47 = "2025-10-19T23:54:27.658355"
    48 = "2025-10-19T23:54:27.659355"
    49 = "2025-10-19T23:54:27.659355"
    50 = "2025-10-19T23:54:27.659355"
    51 = "2025-10-19T23:54:27.659355"
    52 = "2025-10-19T23:54:27.659355"
    53 = "2025-10-19T23:54:27.659355"
    54 = "2025-10-19T23:54:27.659355"
    55 = "2025-10-19T23:54:27.659355"
    56 = "2025-10-19T23:54:27.660353"
    57 = "2025-10-19T23:54:27.660353"
    58 = "2025-10-19T23:54:27.660353"
    59 = "2025-10-19T23:54:27.660353"
    60 = "2025-10-19T23:54:27.660353"
    61 = "2025-10-19T23:54:27.660353"
    62 = "2025-10-19T23:54:27.660353"
    63 = "2025-10-19T23:54:27.660353"
    64 = "2025-10-19T23:54:27.660353"
    65 = "2025-10-19T23:54:27.660353"
    66 = "2025-10-19T23:54:27.661353"
    67 = "2025-10-19T23:54:27.661353"
    68 = "2025-10-19T23:54:27.661353"
    69 = "2025-10-19T23:54:27.661353"
    70 = "2025-10-19T23:54:27.661353"
    71 = "2025-10-19T23:54:27.661353"
    72 = "2025-10-19T23:54:27.661353"
    73 = "2025-10-19T23:54:27.661353"
    74 = "2025-10-19T23:54:27.662353"
    75 = "2025-10-19T23:54:27.662353"
    76 = "2025-10-19T23:54:27.662353"
    77 = "2025-10-19T23:54:27.662353"
    78 = "2025-10-19T23:54:27.662353"
    79 = "2025-10-19T23:54:27.662353"
    80 = "2025-10-19T23:54:27.662353"
    81 = "2025-10-19T23:54:27.662353"
    82 = "2025-10-19T23:54:27.662353"
    83 = "2025-10-19T23:54:27.662353"
    84 = "2025-10-19T23:54:27.663352"
    85 = "2025-10-19T23:54:27.663352"
    86 = "2025-10-19T23:54:27.663352"
    87 = "2025-10-19T23:54:27.664360"
    88 = "2025-10-19T23:54:27.664360"
    89 = "2025-10-19T23:54:27.664360"
    90 = "2025-10-19T23:54:27.664360"
    91 = "2025-10-19T23:54:27.665356"
    92 = "2025-10-19T23:54:27.665356"
    93 = "2025-10-19T23:54:27.665356"
    94 = "2025-10-19T23:54:27.665356"
    95 = "2025-10-19T23:54:27.666357"
    96 = "2025-10-19T23:54:27.666458"
    97 = "2025-10-19T23:54:27.666458"
    98 = "2025-10-19T23:54:27.666994"
    99 = "2025-10-19T23:54:27.666994"
    100 = "2025-10-19T23:54:27.667542"
    101 = "2025-10-19T23:54:27.667542"
    102 = "2025-10-19T23:54:27.667542"
    103 = "2025-10-19T23:54:27.667542"
    104 = "2025-10-19T23:54:27.667542"
    105 = "2025-10-19T23:54:27.668521"
    106 = "2025-10-19T23:54:27.668521"
    107 = "2025-10-19T23:54:27.668521"
    108 = "2025-10-19T23:54:27.669521"
    109 = "2025-10-19T23:54:27.669521"
    110 = "2025-10-19T23:54:27.669521"
    111 = "2025-10-19T23:54:27.669521"
    112 = "2025-10-19T23:54:27.670519"
    113 = "2025-10-19T23:54:27.670519"
    114 = "2025-10-19T23:54:27.670519"
    115 = "2025-10-19T23:54:27.670519"
    116 = "2025-10-19T23:54:27.670519"
    117 = "2025-10-19T23:54:27.672520"
    118 = "2025-10-19T23:54:27.673519"
    119 = "2025-10-19T23:54:27.673519"
    120 = "2025-10-19T23:54:27.673519"
    121 = "2025-10-19T23:54:27.673519"
    122 = "2025-10-19T23:54:27.673519"
    123 = "2025-10-19T23:54:27.673519"
    124 = "2025-10-19T23:54:27.673519"
    125 = "2025-10-19T23:54:27.674518"
    126 = "2025-10-19T23:54:27.674518"
    127 = "2025-10-19T23:54:27.674518"
    128 = "2025-10-19T23:54:27.674518"
    129 = "2025-10-19T23:54:27.674518"
    130 = "2025-10-19T23:54:27.674518"
    131 = "2025-10-19T23:54:27.674518"
    132 = "2025-10-19T23:54:27.674518"
    133 = "2025-10-19T23:54:27.674518"
    134 = "2025-10-19T23:54:27.675518"
    135 = "2025-10-19T23:54:27.675518"
    136 = "2025-10-19T23:54:27.675518"
    137 = "2025-10-19T23:54:27.675518"
    138 = "2025-10-19T23:54:27.675518"
    139 = "2025-10-19T23:54:27.675518"
    140 = "2025-10-19T23:54:27.675518"
    141 = "2025-10-19T23:54:27.675518"
    142 = "2025-10-19T23:54:27.675518"
    143 = "2025-10-19T23:54:27.676518"
    144 = "2025-10-19T23:54:27.676518"
    145 = "2025-10-19T23:54:27.676518"
    146 = "2025-10-19T23:54:27.676518"
    147 = "2025-10-19T23:54:27.676518"
    148 = "2025-10-19T23:54:27.676518"
    149 = "2025-10-19T23:54:27.676518"
    150 = "2025-10-19T23:54:27.676518"
    151 = "2025-10-19T23:54:27.676518"
    152 = "2025-10-19T23:54:27.677518"
    153 = "2025-10-19T23:54:27.677518"
    154 = "2025-10-19T23:54:27.677518"
    155 = "2025-10-19T23:54:27.677518"
    156 = "2025-10-19T23:54:27.677518"
    157 = "2025-10-19T23:54:27.677518"
    158 = "2025-10-19T23:54:27.677518"
    159 = "2025-10-19T23:54:27.677518"
    160 = "2025-10-19T23:54:27.677518"
    161 = "2025-10-19T23:54:27.678518"
    162 = "2025-10-19T23:54:27.678518"
    163 = "2025-10-19T23:54:27.678518"
    164 = "2025-10-19T23:54:27.679521"
    165 = "2025-10-19T23:54:27.679521"
    166 = "2025-10-19T23:54:27.679521"
    167 = "2025-10-19T23:54:27.679521"
    168 = "2025-10-19T23:54:27.680517"
    169 = "2025-10-19T23:54:27.680517"
    170 = "2025-10-19T23:54:27.680517"
    171 = "2025-10-19T23:54:27.680517"
}

foreach ($row in $rowTimestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $rowTimestamps[$row]
}
